$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "update daily + request kak oyko"
# Row 15 (item 7) was an empty template row; fill it in with the next daily
# entry, matching the look/format of the already-filled rows above it
# (row 14) — copy that row's formatting across first, then set the values.

$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B15").Value = 44746
$ws.Range("C15").Value = "- Mining nama nama Toko Watson pada BiVi`n- Mengenal Query data pada BiVi (mandiri), struktur database dan table.`n- mining data data hasil penjulan watson untuk beberapa lokasi`n- EDA data Watson dan BiVi"
$ws.Range("D15").Value = 0.3125
$ws.Range("E15").Value = 0.70833333333333337

# The longer description needs a taller row.
$ws.Rows("15").RowHeight = 104.25

# Leave the selection where the author's cursor ended up.
$ws.Range("H15").Select()
